# repull data, push all data, mean calculation
# Update the "dSF" column (F) values for the rows whose source data was re-pulled.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    6  = 3
    7  = -1
    17 = -1
    24 = 9
    26 = 2
    29 = 4
    31 = 1
    33 = -3
    34 = 1
    42 = 3
    44 = 3
    46 = -8
    47 = 3
    53 = -2
    60 = -3
    61 = -2
    63 = 0
    65 = -6
    66 = 1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
